$d = $word.ActiveDocument

# Locate the paragraph ending with "Faz versionamento (push) dos códigos na nuvem."
# (the second "Faz versionamento" occurrence, inside the "GitHub: repositório remoto..." paragraph)
$findRng = $d.Content
$findRng.Find.Execute("Faz versionamento (push) dos códigos na nuvem.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$targetPara = $findRng.Paragraphs(1)
$paraEnd = $targetPara.Range.End

# Insertion point: right before the paragraph mark (End - 1), so the new
# content is inserted as new paragraphs right after this paragraph without
# disturbing the following paragraph.
$insertRng = $d.Range($paraEnd - 1, $paraEnd - 1)

$xmlFrag = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Tutorial de como acessar o Github em outro computador</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Menu File &gt; accounts  &gt; acessar a opção superior (github.com) &gt; antes de clonar qualquer repositório é necessário ativar o Git &gt;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Menu: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>File/options/git/save &gt; no menu principal</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> (área inicial)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>, escolher qual dos seus repositórios será clonado &gt; escolher o local e selecionar a pasta</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> &gt; trabalhar, comitar o arquivo e fazer push para o Github.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Para sair e não deixar rastros</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Fechar todas os programas, menos o github desktop &gt; Repository &gt; Remove e mandar a pasta para a lixeira &gt;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Menu File &gt; Options &gt; Accounts &gt; Sign out (deslogar seu usuário)  &gt; Menu lateral Git &gt; apagar name e email &gt; Salvar &gt; esvaziar a lixeira.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertRng.InsertXML($xmlFrag)
